# Weekly Logs - add tasks done during the following 2-3 weeks (rows 9, 10, 11)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colB = "• Worked on UI (SignIn, Teach, Settings)"
$colC = "• Worked on the database, fixing errors"
$colD = "• Worked on UI (Home, Learn) Created prototypes and statically populated names."
$colE = "• Worked on UI communicating with the database"

foreach ($r in 9, 10, 11) {
    $ws.Range("B$r").Value = $colB
    $ws.Range("C$r").Value = $colC
    $ws.Range("D$r").Value = $colD
    $ws.Range("E$r").Value = $colE

    # D and E use the wrap-text style (matching B5:E5's formatting)
    $ws.Range("D$r").WrapText = $true
    $ws.Range("E$r").WrapText = $true

    # Rows grow to fit the wrapped text
    $ws.Rows.Item($r).RowHeight = 30
}

# Move the active selection to reflect where editing ended up
$ws.Range("E14").Select()
